# Rename the "congenital" variable-name entry to "misc_long_term" across
# every sheet in the regression name list that still holds the old value.
# Each "variables_NNN" sheet stores its list of variable names in column A
# (A1 = sheet title, A2.. = variable names); the cell holding "congenital"
# moves to "misc_long_term".

$wb = $excel.ActiveWorkbook

$targets = @(
    @{ Sheet = "variables_90";  Cell = "A3" },
    @{ Sheet = "variables_112"; Cell = "A4" },
    @{ Sheet = "variables_120"; Cell = "A3" },
    @{ Sheet = "variables_121"; Cell = "A3" },
    @{ Sheet = "variables_122"; Cell = "A3" },
    @{ Sheet = "variables_123"; Cell = "A3" },
    @{ Sheet = "variables_124"; Cell = "A3" },
    @{ Sheet = "variables_125"; Cell = "A3" },
    @{ Sheet = "variables_126"; Cell = "A3" },
    @{ Sheet = "variables_148"; Cell = "A3" },
    @{ Sheet = "variables_149"; Cell = "A3" },
    @{ Sheet = "variables_150"; Cell = "A3" },
    @{ Sheet = "variables_151"; Cell = "A3" },
    @{ Sheet = "variables_152"; Cell = "A3" },
    @{ Sheet = "variables_153"; Cell = "A3" },
    @{ Sheet = "variables_154"; Cell = "A3" },
    @{ Sheet = "variables_176"; Cell = "A2" },
    @{ Sheet = "variables_177"; Cell = "A2" },
    @{ Sheet = "variables_178"; Cell = "A2" }
)

foreach ($target in $targets) {
    $ws = $wb.Worksheets.Item($target.Sheet)
    $ws.Range($target.Cell).Value = "misc_long_term"
}
